# Apply "repull data, push all data, mean calculation" updates to column F (dSF)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = 4
    7  = -8
    10 = 11
    11 = -1
    16 = 4
    17 = 0
    18 = 0
    21 = 3
    22 = 2
    27 = 0
    28 = 6
    34 = 4
    41 = 1
    42 = -2
    43 = -5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
